# Updated cryptos list on Fri May 12 13:16:03 UTC 2023 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) figures, and fix two rank swaps (HuobiToken/Frax, TheSandbox/Algorand).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds plain-looking numeric strings (e.g. "1.005") that Excel would
# otherwise auto-convert to floating point numbers (losing exact text + introducing FP
# rounding noise). Force the data rows to Text format first so values are stored verbatim,
# matching the workbook's original inlineStr cells, then drop the format back to Normal so
# no stray cell styling is introduced.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.513.22"
$ws.Range("D3").Value = "1.778.79"
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "1.004"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "307.42"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("D7").Value = "0.4319"
$ws.Range("E7").Value = "  +1.88%  "
$ws.Range("D8").Value = "0.3668"
$ws.Range("E8").Value = "  +2.23%  "
$ws.Range("D9").Value = "0.07203"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").Value = "0.8528"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").Value = "20.47"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "1.778.90"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "6.474"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").Value = "5.268"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "0.06896"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "79.55"
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").Value = "0.000008726"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "15.12"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").Value = "26.518.29"
$ws.Range("E21").Value = "  -2.73%  "
$ws.Range("D22").Value = "5.130"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "11.18"
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("D24").Value = "2.001.98"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "152.45"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("D26").Value = "1.874"
$ws.Range("E26").Value = "  -5.37%  "
$ws.Range("D27").Value = "18.19"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").Value = "5.118"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "114.46"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").Value = "1.792"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").Value = "0.08989"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").Value = "0.7300"
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("D33").Value = "1.124"
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("D34").Value = "4.355"
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("B35").Value = "Frax"
$ws.Range("C35").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D35").Value = "1.005"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "2.742"
$ws.Range("E36").Value = "  -6.54%  "
$ws.Range("D37").Value = "1.088"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").Value = "0.05194"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("D39").Value = "0.01895"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "0.1621"
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.4952"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("D42").Value = "2.603"
$ws.Range("E42").Value = "  -6.55%  "
$ws.Range("D43").Value = "6.352"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").Value = "8.073"
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("D45").Value = "105.28"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D47").Value = "1.004"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").Value = "0.4514"
$ws.Range("E48").Value = "  -3.33%  "
$ws.Range("D49").Value = "0.06214"
$ws.Range("E49").Value = "  -3.58%  "
$ws.Range("D50").Value = "1.610"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "1.744"
$ws.Range("E51").Value = "  +3.59%  "

$priceRange.Style = "Normal"
